$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell "D2" "57.126.85"
$ws.Range("E2").Value = "  -8.27%  "

Set-TextCell "D3" "2.855.40"
$ws.Range("E3").Value = "  -8.17%  "

$ws.Range("E4").Value = "  -0.10%  "

Set-TextCell "D5" "544.60"
$ws.Range("E5").Value = "  -8.00%  "

Set-TextCell "D6" "120.32"
$ws.Range("E6").Value = "  -9.28%  "

$ws.Range("E7").Value = "  +0.06%  "

Set-TextCell "D8" "2.846.91"
$ws.Range("E8").Value = "  -8.58%  "

Set-TextCell "D9" "0.484"
$ws.Range("E9").Value = "  -3.97%  "

Set-TextCell "D10" "0.124"
$ws.Range("E10").Value = "  -11.82%  "

Set-TextCell "D11" "4.74"
$ws.Range("E11").Value = "  -10.24%  "

Set-TextCell "D12" "0.422"
$ws.Range("E12").Value = "  -5.06%  "

Set-TextCell "D13" "0.0000209"
$ws.Range("E13").Value = "  -11.97%  "

Set-TextCell "D14" "30.89"
$ws.Range("E14").Value = "  -9.05%  "

$ws.Range("E15").Value = "  -2.21%  "

Set-TextCell "D16" "3.321.85"
$ws.Range("E16").Value = "  -8.76%  "

Set-TextCell "D17" "2.841.48"
$ws.Range("E17").Value = "  -9.37%  "

Set-TextCell "D18" "57.165.09"
$ws.Range("E18").Value = "  -8.83%  "

Set-TextCell "D19" "6.20"
$ws.Range("E19").Value = "  -3.57%  "

Set-TextCell "D20" "406.66"
$ws.Range("E20").Value = "  -10.04%  "

Set-TextCell "D21" "12.56"
$ws.Range("E21").Value = "  -8.11%  "

Set-TextCell "D22" "0.638"
$ws.Range("E22").Value = "  -6.45%  "

Set-TextCell "D23" "6.64"
$ws.Range("E23").Value = "  -11.14%  "

Set-TextCell "D24" "12.31"
$ws.Range("E24").Value = "  -5.01%  "

Set-TextCell "D25" "75.51"
$ws.Range("E25").Value = "  -7.52%  "

Set-TextCell "D26" "1.00"
$ws.Range("E26").Value = "  +0.33%  "

Set-TextCell "D27" "0.998"
$ws.Range("E27").Value = "  -0.49%  "

Set-TextCell "D28" "2.40"
$ws.Range("E28").Value = "  -8.71%  "

$ws.Range("E29").Value = "  -7.93%  "

$ws.Range("E30").Value = "  -7.21%  "

Set-TextCell "D31" "24.36"
$ws.Range("E31").Value = "  -7.93%  "

Set-TextCell "D32" "5.84"
$ws.Range("E32").Value = "  -11.19%  "

Set-TextCell "D33" "0.0915"
$ws.Range("E33").Value = "  -7.51%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D34" "5.28"
$ws.Range("E34").Value = "  -8.60%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D35" "47.97"
$ws.Range("E35").Value = "  -5.42%  "

$ws.Range("B36").Value = "Stacks"
$ws.Range("C36").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D36" "1.96"
$ws.Range("E36").Value = "  -16.36%  "

$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D37" "0.874"
$ws.Range("E37").Value = "  -12.44%  "

Set-TextCell "D38" "8.11"
$ws.Range("E38").Value = "  +1.36%  "

Set-TextCell "D39" "0.0₃0608"
$ws.Range("E39").Value = "  -13.92%  "

Set-TextCell "D40" "0.0334"
$ws.Range("E40").Value = "  -12.38%  "

$ws.Range("E41").Value = "  -7.34%  "

Set-TextCell "D42" "2.571.82"
$ws.Range("E42").Value = "  -6.23%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D44" "2.32"
$ws.Range("E44").Value = "  -9.52%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D45" "342.64"
$ws.Range("E45").Value = "  -10.46%  "

Set-TextCell "D46" "116.57"
$ws.Range("E46").Value = "  -7.28%  "

Set-TextCell "D47" "0.224"
$ws.Range("E47").Value = "  -8.54%  "

Set-TextCell "D48" "0.104"
$ws.Range("E48").Value = "  -5.33%  "

Set-TextCell "D49" "1.89"
$ws.Range("E49").Value = "  -8.20%  "

Set-TextCell "D50" "22.26"
$ws.Range("E50").Value = "  -8.90%  "

Set-TextCell "D51" "1.90"
$ws.Range("E51").Value = "  -10.29%  "
